# Cell.value_without_whitespace fix: add a regression-test cell that holds
# the same text as the existing whitespace-string sample (row 3, col F) but
# WITHOUT the leading newline, so the two can be compared side by side.
#
# F3 already contains shared string "\nwhitespace    string" (index 13).
# G3 gets the trimmed variant "whitespace    string" (a brand-new shared
# string, index 14) with the same wrapped-text style as F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCell = $ws.Range("G3")
$newCell.Value = "whitespace    string"
$newCell.WrapText = $true

# Resize column G the way Excel does after typing a new value into it
# (double-click the column border / Format > AutoFit Column Width).
$ws.Range("G3").EntireColumn.AutoFit() | Out-Null

# Excel leaves the newly-edited cell selected.
$newCell.Select()
